$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '62.741.83'
$ws.Range('E2').Value = '  -2.18%  '
$ws.Range('D3').Value = '3.200.88'
$ws.Range('E3').Value = '  -3.38%  '
$ws.Range('E4').Value = '  +0.01%  '
$ws.Range('D5').Value = '592.21'
$ws.Range('E5').Value = '  -1.40%  '
$ws.Range('D6').Value = '136.33'
$ws.Range('E6').Value = '  -5.29%  '
$ws.Range('E7').Value = '  -0.09%  '
$ws.Range('D8').Value = '3.202.19'
$ws.Range('E8').Value = '  -3.18%  '
$ws.Range('E9').Value = '  -3.03%  '
$ws.Range('E10').Value = '  -3.32%  '
$ws.Range('D11').Value = '5.36'
$ws.Range('E11').Value = '  -2.23%  '
$ws.Range('E12').Value = '  -4.13%  '
$ws.Range('D13').Value = '0.0000239'
$ws.Range('E13').Value = '  -4.34%  '
$ws.Range('D14').Value = '33.58'
$ws.Range('E14').Value = '  -4.11%  '
$ws.Range('D15').Value = '3.730.37'
$ws.Range('E15').Value = '  -3.40%  '
$ws.Range('E16').Value = '  -0.12%  '
$ws.Range('D17').Value = '3.201.82'
$ws.Range('E17').Value = '  -3.44%  '
$ws.Range('D18').Value = '62.805.70'
$ws.Range('E18').Value = '  -2.19%  '
$ws.Range('E19').Value = '  -3.18%  '
$ws.Range('D20').Value = '464.36'
$ws.Range('E20').Value = '  -4.29%  '
$ws.Range('D21').Value = '13.91'
$ws.Range('E21').Value = '  -3.22%  '
$ws.Range('D22').Value = '0.716'
$ws.Range('E22').Value = '  -4.12%  '
$ws.Range('E23').Value = '  -5.11%  '
$ws.Range('D24').Value = '13.42'
$ws.Range('E24').Value = '  -1.44%  '
$ws.Range('D25').Value = '83.87'
$ws.Range('E25').Value = '  -1.15%  '
$ws.Range('E26').Value = '  -0.01%  '
$ws.Range('D27').Value = '2.71'
$ws.Range('E27').Value = '  -2.97%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '1.00'
$ws.Range('D28').ClearFormats()
$ws.Range('E28').Value = '  -0.08%  '
$ws.Range('B29').Value = 'RenderToken'
$ws.Range('C29').Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range('D29').Value = '7.89'
$ws.Range('E29').Value = '  -5.15%  '
$ws.Range('B30').Value = 'NEARProtocol'
$ws.Range('C30').Value = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
$ws.Range('D30').Value = '6.91'
$ws.Range('E30').Value = '  -4.79%  '
$ws.Range('E31').Value = '  -4.12%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '27.60'
$ws.Range('D32').ClearFormats()
$ws.Range('E32').Value = '  -3.42%  '
$ws.Range('E33').Value = '  -4.86%  '
$ws.Range('D34').Value = '2.44'
$ws.Range('E34').Value = '  -5.78%  '
$ws.Range('E35').Value = '  -5.10%  '
$ws.Range('E36').Value = '  -2.73%  '
$ws.Range('D37').Value = '51.64'
$ws.Range('E37').Value = '  -3.65%  '
$ws.Range('D38').Value = '0.0₃0701'
$ws.Range('E38').Value = '  -5.34%  '
$ws.Range('E39').Value = '  -2.87%  '
$ws.Range('D40').Value = '420.15'
$ws.Range('E40').Value = '  -2.85%  '
$ws.Range('D41').Value = '3.002.32'
$ws.Range('E41').Value = '  -0.68%  '
$ws.Range('E42').Value = '  +4.17%  '
$ws.Range('E43').Value = '  -4.42%  '
$ws.Range('D44').Value = '2.63'
$ws.Range('E44').Value = '  -6.06%  '
$ws.Range('D45').Value = '0.255'
$ws.Range('E45').Value = '  -6.28%  '
$ws.Range('E46').Value = '  -4.86%  '
$ws.Range('D48').Value = '35.38'
$ws.Range('E48').Value = '  -0.14%  '
$ws.Range('D49').Value = '25.67'
$ws.Range('E49').Value = '  -2.84%  '
$ws.Range('D50').Value = '125.25'
$ws.Range('E50').Value = '  +1.40%  '
$ws.Range('E51').Value = '  -2.96%  '
